# Daily attendance processing - 2025-11-04 07:44:36
# Normalize the "Recorded By" (column G) entries: when the system-generated
# "System" tag is listed first, move it to the end of the comma-separated
# list of recorders instead (e.g. "System, user@x.com" -> "user@x.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$prefix = "System, "

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Text

    if ($current.StartsWith($prefix)) {
        $remainder = $current.Substring($prefix.Length)
        $cell.Value = $remainder + ", System"
    }
}
